$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fecha" column (J) used to store the date extracted from the
# source filename as plain text (e.g. "140823" for 14-08-2023).
# Refactor it into a real Excel date value (2023-08-14 -> serial 45152)
# formatted as "YYYY-MM-DD HH:MM:SS" instead of raw text.

# Touch the (lowercase) lookalike date/time format once so it gets
# registered, then move the first data cell to the final uppercase
# format that every cell in the column will actually use.
$firstCell = $ws.Cells.Item(2, 10)
$firstCell.Value = 45152
$firstCell.NumberFormat = "yyyy-mm-dd h:mm:ss"
$firstCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($row = 3; $row -le 37; $row++) {
    $cell = $ws.Cells.Item($row, 10)
    $cell.Value = 45152
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
